$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "51.742.10"
$ws.Range("E2").Value = "  -0.45%  "
Set-TextValue "D3" "2.953.55"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "352.83"
$ws.Range("E5").Value = "  -1.39%  "
Set-TextValue "D6" "105.55"
$ws.Range("E6").Value = "  -4.63%  "
Set-TextValue "D7" "0.549"
$ws.Range("E7").Value = "  -3.44%  "
Set-TextValue "D9" "0.598"
$ws.Range("E9").Value = "  -5.43%  "
Set-TextValue "D10" "37.44"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("E11").Value = "  +2.16%  "
Set-TextValue "D12" "0.0846"
$ws.Range("E12").Value = "  -3.86%  "
Set-TextValue "D13" "18.81"
$ws.Range("E13").Value = "  -4.79%  "
Set-TextValue "D14" "3.433.06"
$ws.Range("E14").Value = "  +1.09%  "
Set-TextValue "D15" "7.45"
$ws.Range("E15").Value = "  -5.92%  "
Set-TextValue "D16" "2.944.75"
$ws.Range("E16").Value = "  +0.29%  "
Set-TextValue "D17" "0.983"
$ws.Range("E17").Value = "  -0.65%  "
Set-TextValue "D18" "51.702.96"
$ws.Range("E18").Value = "  -0.57%  "
Set-TextValue "D19" "3.31"
$ws.Range("E19").Value = "  -0.91%  "
Set-TextValue "D20" "7.31"
$ws.Range("E20").Value = "  -3.70%  "
Set-TextValue "D21" "13.26"
$ws.Range("E21").Value = "  -5.69%  "
Set-TextValue "D22" "0.0₃0957"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D23" "265.58"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "68.84"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("E25").Value = "  -5.05%  "
$ws.Range("E26").Value = "  -7.62%  "
Set-TextValue "D27" "26.53"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D29" "7.23"
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.109"
$ws.Range("E30").Value = "  +2.75%  "
Set-TextValue "D31" "6.23"
$ws.Range("E31").Value = "  +2.20%  "
Set-TextValue "D32" "10.02"
$ws.Range("E32").Value = "  -5.98%  "
Set-TextValue "D33" "2.16"
$ws.Range("E33").Value = "  -5.32%  "
Set-TextValue "D34" "35.34"
$ws.Range("E34").Value = "  -7.65%  "
Set-TextValue "D35" "50.76"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.0425"
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D40" "17.24"
$ws.Range("E40").Value = "  -6.86%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D41" "1.91"
$ws.Range("E41").Value = "  -5.42%  "
Set-TextValue "D42" "0.115"
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D43" "123.34"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "22.60"
$ws.Range("E44").Value = "  -1.84%  "
Set-TextValue "D45" "2.15"
$ws.Range("E45").Value = "  -0.66%  "
Set-TextValue "D46" "2.104.38"
$ws.Range("E46").Value = "  -1.61%  "
Set-TextValue "D47" "3.24"
$ws.Range("E47").Value = "  -6.74%  "
$ws.Range("E48").Value = "  -7.61%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D49" "3.262.43"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D50" "0.237"
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-TextValue "D51" "0.0322"
$ws.Range("E51").Value = "  -3.37%  "
